# Auto-generated Excel COM-interop script
# Applies scheduled-runner price/profit refresh to ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 24717.143
$ws.Range("J75").Value = 24717.143
$ws.Range("L75").Value = 24717.143
$ws.Range("N75").Value = -26589.143

$ws.Range("H78").Value = 24717.143
$ws.Range("J78").Value = 24717.143
$ws.Range("L78").Value = 74151.429
$ws.Range("N78").Value = -83511.429

$ws.Range("H88").Value = 3165.3125
$ws.Range("I88").Value = 699.5
$ws.Range("J88").Value = 10562.75
$ws.Range("K88").Value = 699.5
$ws.Range("L88").Value = 10562.75
$ws.Range("M88").Value = -293.5
$ws.Range("N88").Value = -11374.75

$ws.Range("H91").Value = 3165.3125
$ws.Range("I91").Value = 699.5
$ws.Range("J91").Value = 10562.75
$ws.Range("K91").Value = 699.5
$ws.Range("L91").Value = 10562.75
$ws.Range("M91").Value = 704.5
$ws.Range("N91").Value = -13370.75

$ws.Range("H118").Value = 750
$ws.Range("I118").Value = 392.85715
$ws.Range("J118").Value = 2000
$ws.Range("K118").Value = 1178.57145
$ws.Range("L118").Value = 6000
$ws.Range("M118").Value = 478.4285500000001
$ws.Range("N118").Value = -9314

$ws.Range("H137").Value = 61477.883
$ws.Range("I137").Value = 2820.5
$ws.Range("K137").Value = 8461.5
$ws.Range("M137").Value = -5911.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3999
$ws.Range("I45").Value = 3999
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 3999
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -3622
$ws.Range("N45").ClearContents()

$ws.Range("H88").Value = 2130.9285
$ws.Range("I88").Value = 2115.111
$ws.Range("J88").Value = 2159.4
$ws.Range("K88").Value = 2115.111
$ws.Range("L88").Value = 2159.4
$ws.Range("M88").Value = -1709.111
$ws.Range("N88").Value = -2971.4

$ws.Range("H91").Value = 2130.9285
$ws.Range("I91").Value = 2115.111
$ws.Range("J91").Value = 2159.4
$ws.Range("K91").Value = 2115.111
$ws.Range("L91").Value = 2159.4
$ws.Range("M91").Value = -711.1109999999999
$ws.Range("N91").Value = -4967.4

$ws.Range("H101").Value = 36666.4
$ws.Range("J101").Value = 36666.4
$ws.Range("L101").Value = 36666.4
$ws.Range("N101").Value = -43156.4

$ws.Range("H102").Value = 2227.5
$ws.Range("I102").Value = 1470
$ws.Range("J102").Value = 4500
$ws.Range("K102").Value = 1470
$ws.Range("L102").Value = 4500
$ws.Range("M102").Value = 152
$ws.Range("N102").Value = -7744

$ws.Range("H123").Value = 40609.668
$ws.Range("J123").Value = 40609.668
$ws.Range("L123").Value = 40609.668
$ws.Range("N123").Value = -50409.668

$ws.Range("H132").Value = 3114.9333
$ws.Range("I132").Value = 1840.5714
$ws.Range("J132").Value = 4230
$ws.Range("K132").Value = 5521.7142
$ws.Range("L132").Value = 12690
$ws.Range("M132").Value = -2991.7142
$ws.Range("N132").Value = -17750

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2201.4443
$ws.Range("I86").Value = 1900
$ws.Range("J86").Value = 2578.25
$ws.Range("K86").Value = 1900
$ws.Range("L86").Value = 2578.25
$ws.Range("M86").Value = -777
$ws.Range("N86").Value = -4824.25

$ws.Range("H89").Value = 2201.4443
$ws.Range("I89").Value = 1900
$ws.Range("J89").Value = 2578.25
$ws.Range("K89").Value = 9500
$ws.Range("L89").Value = 12891.25
$ws.Range("M89").Value = -3884
$ws.Range("N89").Value = -24123.25

$ws.Range("H99").Value = 2839.476
$ws.Range("I99").Value = 2151.5833
$ws.Range("J99").Value = 3756.6667
$ws.Range("K99").Value = 2151.5833
$ws.Range("L99").Value = 3756.6667
$ws.Range("M99").Value = -653.5832999999998
$ws.Range("N99").Value = -6752.6667

$ws.Range("H134").Value = 2434
$ws.Range("I134").Value = 2504
$ws.Range("J134").Value = 2014
$ws.Range("K134").Value = 7512
$ws.Range("L134").Value = 6042
$ws.Range("M134").Value = -4977
$ws.Range("N134").Value = -11112

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2450
$ws.Range("I16").Value = 1450
$ws.Range("J16").Value = 6450
$ws.Range("K16").Value = 1450
$ws.Range("L16").Value = 6450
$ws.Range("M16").Value = -1163
$ws.Range("N16").Value = -7024

$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()

$ws.Range("H74").Value = 26471.428
$ws.Range("J74").Value = 26875
$ws.Range("L74").Value = 26875
$ws.Range("N74").Value = -28623

$ws.Range("H77").Value = 26471.428
$ws.Range("J77").Value = 26875
$ws.Range("L77").Value = 80625
$ws.Range("N77").Value = -89361

$ws.Range("H99").Value = 2700
$ws.Range("I99").Value = 2700
$ws.Range("K99").Value = 2700
$ws.Range("M99").Value = -1202

$ws.Range("H107").Value = 645.4815
$ws.Range("I107").Value = 599.3333
$ws.Range("J107").Value = 703.1667
$ws.Range("K107").Value = 599.3333
$ws.Range("L107").Value = 703.1667
$ws.Range("M107").Value = 1320.6667
$ws.Range("N107").Value = -4543.1667

$ws.Range("H113").Value = 2450
$ws.Range("I113").Value = 1450
$ws.Range("J113").Value = 6450
$ws.Range("K113").Value = 1450
$ws.Range("L113").Value = 6450
$ws.Range("M113").Value = 720
$ws.Range("N113").Value = -10790

$ws.Range("H126").Value = 2700
$ws.Range("I126").Value = 2700
$ws.Range("K126").Value = 8100
$ws.Range("M126").Value = -5630

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 944
$ws.Range("I5").Value = 1006.7273
$ws.Range("J5").Value = 881.2727
$ws.Range("K5").Value = 3020.1819
$ws.Range("L5").Value = 2643.8181
$ws.Range("M5").Value = -2908.1819
$ws.Range("N5").Value = -2867.8181

$ws.Range("H131").Value = 804.28845
$ws.Range("J131").Value = 952.5143
$ws.Range("L131").Value = 2857.5429
$ws.Range("N131").Value = -12937.5429

$ws.Range("H135").Value = 944
$ws.Range("I135").Value = 1006.7273
$ws.Range("J135").Value = 881.2727
$ws.Range("K135").Value = 9060.545700000001
$ws.Range("L135").Value = 7931.454299999999
$ws.Range("M135").Value = -6525.545700000001
$ws.Range("N135").Value = -13001.4543

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 173.57143
$ws.Range("I107").Value = 193.3125
$ws.Range("J107").Value = 110.4
$ws.Range("K107").Value = 193.3125
$ws.Range("L107").Value = 110.4
$ws.Range("M107").Value = 1726.6875
$ws.Range("N107").Value = -3950.4

$ws.Range("H132").Value = 6568.643
$ws.Range("I132").Value = 7498.25
$ws.Range("J132").Value = 5329.1665
$ws.Range("K132").Value = 22494.75
$ws.Range("L132").Value = 15987.4995
$ws.Range("M132").Value = -19964.75
$ws.Range("N132").Value = -21047.4995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 9041.583000000001
$ws.Range("I61").Value = 25449.75
$ws.Range("K61").Value = 25449.75
$ws.Range("M61").Value = -25247.75

$ws.Range("H82").Value = 1557.4286
$ws.Range("I82").Value = 1634
$ws.Range("J82").Value = 1500
$ws.Range("K82").Value = 1634
$ws.Range("L82").Value = 1500
$ws.Range("M82").Value = -1273
$ws.Range("N82").Value = -2222

$ws.Range("H85").Value = 1557.4286
$ws.Range("I85").Value = 1634
$ws.Range("J85").Value = 1500
$ws.Range("K85").Value = 1634
$ws.Range("L85").Value = 1500
$ws.Range("M85").Value = -386
$ws.Range("N85").Value = -3996

$ws.Range("H113").Value = 9041.583000000001
$ws.Range("I113").Value = 25449.75
$ws.Range("K113").Value = 25449.75
$ws.Range("M113").Value = -23279.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 336.86667
$ws.Range("I113").Value = 348
$ws.Range("J113").Value = 320.16666
$ws.Range("K113").Value = 1044
$ws.Range("L113").Value = 960.4999799999999
$ws.Range("M113").Value = 1126
$ws.Range("N113").Value = -5300.49998

$ws.Range("H114").Value = 60000
$ws.Range("J114").Value = 60000
$ws.Range("L114").Value = 60000
$ws.Range("N114").Value = -68678

$ws.Range("H132").Value = 1948.711
$ws.Range("I132").Value = 1453.7576
$ws.Range("J132").Value = 3309.8333
$ws.Range("K132").Value = 4361.2728
$ws.Range("L132").Value = 9929.499899999999
$ws.Range("M132").Value = -1831.2728
$ws.Range("N132").Value = -14989.4999

